# Applies the post-edit cell values for Proyecto1_Sistema_Sismico/BD.xlsx, sheet "Hoja1".
# Every data cell in this sheet is stored as text (shared string) in the source workbook,
# including values that look like numbers ("8.0") or dates ("24/08/2028"). Each write below
# forces Text format ("@") before assigning the value so Excel does not silently coerce the
# string into a number/date, then resets the cell back to the default "Normal" style so no
# stray per-cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$edits = @(
    @("A2", "24/08/2028"),
    @("B2", "23:07:45"),
    @("C2", "8.0"),
    @("D2", "SUBDUCCION_PLACA"),
    @("E2", "detalle necesario"),
    @("F2", "5.0"),
    @("G2", "4.0"),
    @("H2", "67.9"),
    @("J2", "cambiossss"),
    @("A3", "24/01/2018"),
    @("B3", "20:02:45"),
    @("D3", "CHOQUE_PLACAS"),
    @("E3", "detalle"),
    @("F3", "89.0"),
    @("G3", "9.7"),
    @("H3", "6.0"),
    @("I3", "SAN_JOSE"),
    @("J3", "Hice cambio"),
    @("I10", "ALAJUELA"),
    @("I22", "CARTAGO"),
    @("A28", "24/07/2016"),
    @("J28", "Probar si sirve"),
    @("A30", "24/08/2028"),
    @("B30", "23:07:45"),
    @("C30", "8.0"),
    @("D30", "SUBDUCCION_PLACA"),
    @("E30", "detalle necesario"),
    @("F30", "5.0"),
    @("G30", "4.0"),
    @("H30", "67.9"),
    @("I30", "SAN_JOSE"),
    @("J30", "cambio"),
    @("A31", "24/01/2018"),
    @("B31", "20:02:45"),
    @("C31", "4.0"),
    @("D31", "CHOQUE_PLACAS"),
    @("E31", "detalle"),
    @("F31", "89.0"),
    @("G31", "9.7"),
    @("H31", "6.0"),
    @("I31", "SAN_JOSE"),
    @("J31", "vanes"),
    @("A32", "24/02/2020"),
    @("B32", "02:25:19"),
    @("C32", "0.0"),
    @("D32", "SUBDUCCION_PLACA"),
    @("E32", " Ingrese los detalles"),
    @("F32", "0.0"),
    @("G32", "0.0"),
    @("H32", "0.0"),
    @("I32", "SAN_JOSE"),
    @("J32", " Más detalles"),
    @("A33", "24/12/2016"),
    @("B33", "22:02:20"),
    @("C33", "0.0"),
    @("D33", "TECTONICO_SUBDUCCION"),
    @("E33", "detalle"),
    @("F33", "89.8"),
    @("G33", "5.0"),
    @("H33", "6.0"),
    @("I33", "ALAJUELA"),
    @("J33", "Descripcion")
)

foreach ($edit in $edits) {
    $address = $edit[0]
    $text = $edit[1]
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}
